# Update cryptos list (simulated GitHub Actions refresh of coinranking data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: force "Price" column cells to text format before assigning,
# so numeric-looking strings (e.g. "568.64", "1.00") are preserved exactly
# instead of being coerced into floating point numbers by Excel.

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.124.17"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.546.94"
$ws.Range("E3").Value = "  +3.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.64"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.45"
$ws.Range("E6").Value = "  +3.13%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -0.42%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.546.53"
$ws.Range("E9").Value = "  +3.24%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.14%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  -1.61%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.75%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +0.28%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.49"
$ws.Range("E14").Value = "  +4.25%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.006.81"
$ws.Range("E15").Value = "  +3.35%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.114.93"

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +2.08%  "

# Row 18 - WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.546.29"
$ws.Range("E18").Value = "  +3.20%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +2.22%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.83"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21 - Polkadot
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.30"
$ws.Range("E21").Value = "  +0.86%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.35%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.13%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.38"
$ws.Range("E24").Value = "  -0.14%  "

# Rows 25 & 26 swap places: Fetch.AI <-> Kaspa, with new values
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  -1.97%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.63"
$ws.Range("E26").Value = "  +9.08%  "

# Row 27 - SuiNetwork
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.51"
$ws.Range("E27").Value = "  +11.46%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +5.32%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30 - Aptos
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("E30").Value = "  +7.71%  "

# Row 31 - PEPE
$ws.Range("E31").Value = "  +2.72%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.78%  "

# Row 33 - Monero
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.44"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  +3.85%  "

# Row 35 - Bittensor
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "412.35"
$ws.Range("E35").Value = "  +12.76%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  +0.40%  "

# Row 37 - EthereumClassic
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.04"
$ws.Range("E37").Value = "  +0.97%  "

# Row 38 - NEARProtocol
$ws.Range("E38").Value = "  -0.09%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  -0.02%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +4.30%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42 - OKB
$ws.Range("E42").Value = "  -2.95%  "

# Row 43 - Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "153.43"
$ws.Range("E43").Value = "  +2.13%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +2.13%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.96"
$ws.Range("E45").Value = "  +2.13%  "

# Row 46 - Mantle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.605"
$ws.Range("E46").Value = "  +1.19%  "

# Row 47 - Stellar
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0963"
$ws.Range("E47").Value = "  +0.39%  "

# Row 48 - VeChain
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0239"
$ws.Range("E48").Value = "  +5.86%  "

# Row 49 - Hedera
$ws.Range("E49").Value = "  +1.28%  "

# Row 50 - EnergySwap
$ws.Range("E50").Value = "  +2.04%  "

# Row 51 - dogwifhat
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.78"
$ws.Range("E51").Value = "  +2.46%  "
